$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 'MSG: None

MSG: The decision about which movie to show on Friday was not made, leading to no selection being finalized.
'
$ws.Range("D2").Value = 'no_decision, '
$ws.Range("C3").Value = 'MSG: None

MSG: The decision has been made that no movie will be shown on Friday.
'
$ws.Range("D3").Value = 'no_decision, '
$ws.Range("C4").Value = 'MSG: None

MSG: No decision was made about the movie to be shown on Friday.
'
$ws.Range("D4").Value = 'no_decision, '
$ws.Range("C5").Value = 'MSG: None

MSG: The decision was made to not select any movie for showing on Friday, as the committee could not reach a consensus.
'
$ws.Range("D5").Value = 'no_decision, '
$ws.Range("C6").Value = 'MSG: None

MSG: The decision has been made to acquire the rights for "Barbie."
'
$ws.Range("D6").Value = 'Barbie_was_selected, '
$ws.Range("C7").Value = 'MSG: None

MSG: The decision regarding the movie to show on Friday has been recorded as no decision.
'
$ws.Range("D7").Value = 'no_decision, '
$ws.Range("C8").Value = 'MSG: None

MSG: The decision-making process did not reach a consensus, so no movie has been selected for Friday.
'
$ws.Range("D8").Value = 'no_decision, '
$ws.Range("C9").Value = 'MSG: None

MSG: The decision to show a movie on Friday was not reached during the meeting.
'
$ws.Range("D9").Value = 'no_decision, '
$ws.Range("C10").Value = 'MSG: None

MSG: The decision regarding the movie selection was not finalized, so there will be no acquisition of rights for a movie at this time.
'
$ws.Range("D10").Value = 'no_decision, '
$ws.Range("C11").Value = 'MSG: None

MSG: I have recorded the decision to acquire the rights for the movie "Barbie."
'
$ws.Range("D11").Value = 'Barbie_was_selected, '
$ws.Range("C12").Value = 'MSG: None

MSG: The decision process concluded without reaching an agreement on the movie for Friday.
'
$ws.Range("D12").Value = 'no_decision, '
$ws.Range("C13").Value = 'MSG: None

MSG: The committee has not reached a decision about which movie to show on Friday.
'
$ws.Range("D13").Value = 'no_decision, '
$ws.Range("C14").Value = 'MSG: None

MSG: The decision about Friday''s movie has not been made.
'
$ws.Range("D14").Value = 'no_decision, '
$ws.Range("C15").Value = 'MSG: None

MSG: The decision has been recorded as no decision was reached regarding the movie for Friday.
'
$ws.Range("D15").Value = 'no_decision, '
$ws.Range("C16").Value = 'MSG: None

MSG: The decision regarding the movie to be shown on Friday was not made.
'
$ws.Range("D16").Value = 'no_decision, '
$ws.Range("C17").Value = 'MSG: None

MSG: The decision has been recorded, indicating that no agreement was reached regarding the movie for Friday.
'
$ws.Range("D17").Value = 'no_decision, '
$ws.Range("C18").Value = 'MSG: None

MSG: The decision has been recorded, indicating that no movie was chosen for Friday.
'
$ws.Range("D18").Value = 'no_decision, '
$ws.Range("C19").Value = 'MSG: None

MSG: The decision has been made, and no movie was selected for Friday.
'
$ws.Range("D19").Value = 'no_decision, '
$ws.Range("C20").Value = 'MSG: None

MSG: The decision on which movie to show on Friday has not been made.
'
$ws.Range("D20").Value = 'no_decision, '
$ws.Range("C21").Value = 'MSG: None

MSG: It appears that the committee has not reached a decision about which movie to show on Friday. Therefore, I am submitting the no_decision function.
'
$ws.Range("D21").Value = 'no_decision, '
$ws.Range("C22").Value = 'MSG: None

MSG: The conversation did not result in a decision about which movie to show on Friday.
'
$ws.Range("D22").Value = 'no_decision, '
$ws.Range("C23").Value = 'MSG: None

MSG: The decision about the movie for Friday was not made.
'
$ws.Range("D23").Value = 'no_decision, '
$ws.Range("C24").Value = 'MSG: None

MSG: The decision has been recorded as no decision about Friday''s movie was made.
'
$ws.Range("D24").Value = 'no_decision, '
$ws.Range("C25").Value = 'MSG: None

MSG: The decision has been recorded as no choice of a movie can be made at this time.
'
$ws.Range("D25").Value = 'no_decision, '
$ws.Range("C26").Value = 'MSG: None

MSG: The decision about which movie to show on Friday was not made.
'
$ws.Range("D26").Value = 'no_decision, '
$ws.Range("C27").Value = 'MSG: None

MSG: The decision has been recorded as no decision about Friday''s movie can be made.
'
$ws.Range("D27").Value = 'no_decision, '
$ws.Range("C28").Value = 'MSG: None

MSG: The decision regarding the movie to be shown on Friday resulted in no conclusion, indicating that no selection was made during this discussion.
'
$ws.Range("D28").Value = 'no_decision, '
$ws.Range("C29").Value = 'MSG: None

MSG: The decision made by the committee was to select "Barbie" for the movie to be shown on Friday.
'
$ws.Range("D29").Value = 'Barbie_was_selected, '
$ws.Range("C30").Value = 'MSG: None

MSG: The decision regarding which movie to show on Friday was not made during the discussion.
'
$ws.Range("D30").Value = 'no_decision, '
$ws.Range("C31").Value = 'MSG: None

MSG: The decision has been made that no movie will be acquired for Friday.
'
$ws.Range("D31").Value = 'no_decision, '
$ws.Range("C32").Value = 'MSG: None

MSG: The decision about which movie to show on Friday has been deemed as no decision.
'
$ws.Range("D32").Value = 'no_decision, '
$ws.Range("C33").Value = 'MSG: None

MSG: The decision regarding the movie for Friday was not reached, so no movie will be acquired.
'
$ws.Range("D33").Value = 'no_decision, '
$ws.Range("C34").Value = 'MSG: None

MSG: The decision to acquire the rights for both movies has been successfully recorded.
'
$ws.Range("D34").Value = 'both_movies, '
$ws.Range("C35").Value = 'MSG: None

MSG: The decision was made, and no specific movie was chosen to be shown on Friday.
'
$ws.Range("D35").Value = 'no_decision, '
$ws.Range("C36").Value = 'MSG: None

MSG: The committee did not reach a decision regarding the movie for Friday, and thus, the function for no decision has been recorded.
'
$ws.Range("D36").Value = 'no_decision, '
$ws.Range("C37").Value = 'MSG: None

MSG: The decision resulted in no agreement on which movie to show on Friday.
'
$ws.Range("D37").Value = 'no_decision, '
$ws.Range("C38").Value = 'MSG: None

MSG: I have recorded the decision as "no_decision" since the committee did not arrive at a conclusion about which movie to show on Friday.
'
$ws.Range("D38").Value = 'no_decision, '
$ws.Range("C39").Value = 'MSG: None

MSG: The decision has been recorded as no decision being made regarding the movie for Friday.
'
$ws.Range("D39").Value = 'no_decision, '
$ws.Range("C40").Value = 'MSG: None

MSG: The committee did not reach a decision regarding the movie to be shown on Friday, so I have recorded this as no decision being made.
'
$ws.Range("D40").Value = 'no_decision, '
$ws.Range("C41").Value = 'MSG: None

MSG: It seems I made an error in interpreting the information provided. Since the committee decided on "Barbie" as the movie for Friday''s showing, I will call the appropriate function for that decision. 
```python
Barbie_was_selected()
```
'
$ws.Range("D41").Value = 'both_movies, '
$ws.Range("C42").Value = 'MSG: None

MSG: The decision has been made to not select a movie for Friday.
'
$ws.Range("D42").Value = 'no_decision, '
$ws.Range("C43").Value = 'MSG: None

MSG: The decision-making process did not result in a selection for Friday''s movie, so the outcome is that no decision was made.
'
$ws.Range("D43").Value = 'no_decision, '
$ws.Range("C44").Value = 'MSG: None

MSG: The decision process resulted in no agreement on a movie for Friday.
'
$ws.Range("D44").Value = 'no_decision, '
$ws.Range("C45").Value = 'MSG: None

MSG: No movie was selected in this meeting.
'
$ws.Range("D45").Value = 'no_decision, '
$ws.Range("C46").Value = 'MSG: None

MSG: The decision has been recorded, indicating that no movie was selected for Friday.
'
$ws.Range("D46").Value = 'no_decision, '
$ws.Range("C47").Value = 'MSG: None

MSG: The decision regarding which movie to show on Friday was not finalized, resulting in no decision being made.
'
$ws.Range("D47").Value = 'no_decision, '
$ws.Range("C48").Value = 'MSG: None

MSG: No decision was made regarding the movie to be shown on Friday.
'
$ws.Range("D48").Value = 'no_decision, '
$ws.Range("C49").Value = 'MSG: None

MSG: The decision concluded with no clear choice made for Friday''s movie.
'
$ws.Range("D49").Value = 'no_decision, '
$ws.Range("C50").Value = 'MSG: None

MSG: I have successfully recorded the decision to acquire the rights for the movie "Barbie."
'
$ws.Range("D50").Value = 'Barbie_was_selected, '
$ws.Range("C51").Value = 'MSG: None

MSG: The decision-making committee did not reach a conclusion about which movie to show on Friday, so the function for no decision has been executed.
'
$ws.Range("D51").Value = 'no_decision, '
$ws.Range("C52").Value = 'MSG: None

MSG: The decision process concluded without selecting a movie for Friday.
'
$ws.Range("D52").Value = 'no_decision, '
$ws.Range("C53").Value = 'MSG: None

MSG: The rights to both movies have been acquired.
'
$ws.Range("D53").Value = 'both_movies, '
$ws.Range("C54").Value = 'MSG: None

MSG: The decision process concluded without an agreement on which movie to show on Friday.
'
$ws.Range("D54").Value = 'no_decision, '
$ws.Range("C55").Value = 'MSG: None

MSG: The decision has been recorded, and it appears that no movie was selected by the committee for Friday.
'
$ws.Range("D55").Value = 'no_decision, '
$ws.Range("C56").Value = 'MSG: None

MSG: The committee''s discussion concluded without reaching a decision on the movie for Friday, so no acquisition will take place.
'
$ws.Range("D56").Value = 'no_decision, '
$ws.Range("C57").Value = 'MSG: None

MSG: The decision-making process concluded without a definitive agreement on the movie to show on Friday. Therefore, no further action will be taken regarding acquiring movie rights.
'
$ws.Range("D57").Value = 'no_decision, '
$ws.Range("C58").Value = 'MSG: None

MSG: The decision to acquire the rights to "Barbie" has been recorded successfully.
'
$ws.Range("D58").Value = 'Barbie_was_selected, '
$ws.Range("C59").Value = 'MSG: None

MSG: The decision process has concluded without a choice for Friday’s movie.
'
$ws.Range("D59").Value = 'no_decision, '
$ws.Range("C60").Value = 'MSG: None

MSG: I have recorded the decision as "no decision" regarding which movie to show on Friday.
'
$ws.Range("D60").Value = 'no_decision, '
